$wb = $excel.ActiveWorkbook

# ---- PIR sheet ----
$ws = $wb.Worksheets.Item("PIR")
$rows = @(
    @("2026-01-28", "16:59:54", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:59:57", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:59:59", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:02", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:08", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:11", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:16", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:21", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:27", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:31", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:36", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:41", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:47", "17:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "17:00:52", "17:00", "Bathroom", "No Motion", "Inactive")
)
$r = 250
foreach ($row in $rows) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $row[$c - 1]
        if ($c -eq 1) {
            $cell.NumberFormat = "@"
        } elseif ($v -match "^\d+(\.\d+)?%$") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $v
    }
    $r++
}

# ---- Humidity sheet ----
$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
    @("2026-01-28", "16:59:52", "16:00", "Bathroom", "87.9%", "Active"),
    @("2026-01-28", "16:59:55", "16:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-28", "17:00:00", "17:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-28", "17:00:03", "17:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-28", "17:00:06", "17:00", "Bathroom", "86.9%", "Active"),
    @("2026-01-28", "17:00:09", "17:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-28", "17:00:13", "17:00", "Bathroom", "86.9%", "Active"),
    @("2026-01-28", "17:00:22", "17:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-28", "17:00:26", "17:00", "Bathroom", "86.9%", "Active"),
    @("2026-01-28", "17:00:34", "17:00", "Bathroom", "86.9%", "Active"),
    @("2026-01-28", "17:00:38", "17:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-28", "17:00:43", "17:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-28", "17:00:46", "17:00", "Bathroom", "86.8%", "Active")
)
$r = 246
foreach ($row in $rows) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $row[$c - 1]
        if ($c -eq 1) {
            $cell.NumberFormat = "@"
        } elseif ($v -match "^\d+(\.\d+)?%$") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $v
    }
    $r++
}

# ---- Temperature sheet ----
$ws = $wb.Worksheets.Item("Temperature")
$rows = @(
    @("2026-01-28", "16:59:53", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:59:56", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:01", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:04", "17:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "17:00:07", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:10", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:14", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:23", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:26", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:34", "17:00", "Bathroom", "22.9C", "Active"),
    @("2026-01-28", "17:00:38", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:44", "17:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "17:00:47", "17:00", "Bathroom", "22.8C", "Active")
)
$r = 245
foreach ($row in $rows) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $row[$c - 1]
        if ($c -eq 1) {
            $cell.NumberFormat = "@"
        } elseif ($v -match "^\d+(\.\d+)?%$") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $v
    }
    $r++
}

# ---- mmWave sheet ----
$ws = $wb.Worksheets.Item("mmWave")
$rows = @(
    @("2026-01-28", "16:59:55", "16:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "16:59:57", "16:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "16:59:58", "16:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:00", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:02", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:05", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:05", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:08", "17:00", "Living Room", "PRESENCE", "Active"),
    @("2026-01-28", "17:00:12", "17:00", "Living Room", "PRESENCE", "Active"),
    @("2026-01-28", "17:00:15", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:17", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:21", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:23", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:28", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:29", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:32", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:35", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:39", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:42", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:44", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:48", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:00:50", "17:00", "Living Room", "NO_PRESENCE", "Inactive")
)
$r = 46
foreach ($row in $rows) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $row[$c - 1]
        if ($c -eq 1) {
            $cell.NumberFormat = "@"
        } elseif ($v -match "^\d+(\.\d+)?%$") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $v
    }
    $r++
}
